$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.265550255775452
$ws.Range("B1").Value = 1.011569857597351
$ws.Range("C1").Value = 2.878094673156738
$ws.Range("D1").Value = 5.149636745452881
$ws.Range("E1").Value = 0.9347149729728699
